$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds "Fitness" values for rows 2..252 (Generation 0..250).
# Rows 2..206  (Generation 0..204)   : 7573 -> 7310
# Rows 207..252 (Generation 205..250): 7573 -> 7293

$ws.Range("C2:C206").Value = 7310
$ws.Range("C207:C252").Value = 7293
